$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "<do>"
$ws.Range("C2").Value = 42

# Row 4
$ws.Range("B4").Value = "<down>"

# Row 5
$ws.Range("B5").Value = "<you>"
$ws.Range("C5").Value = 31

# Row 6
$ws.Range("C6").Value = 30

# Row 7
$ws.Range("B7").Value = "<line>"
$ws.Range("C7").Value = 32

# Row 8
$ws.Range("C8").Value = 29

# Row 9
$ws.Range("C9").Value = 36

# Row 10
$ws.Range("C10").Value = 35

# Row 11
$ws.Range("C11").Value = 30

# Row 12
$ws.Range("C12").Value = 28

# Row 13
$ws.Range("C13").Value = 37

# Row 14
$ws.Range("B14").Value = "<the>"
$ws.Range("C14").Value = 38

# Row 15
$ws.Range("B15").Value = "<at>"
$ws.Range("C15").Value = 29

# Row 16
$ws.Range("B16").Value = "<paste>"
$ws.Range("C16").Value = 35

# Row 17
$ws.Range("B17").Value = "<sere>"
$ws.Range("C17").Value = 32

# Row 18
$ws.Range("C18").Value = 32
